$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 300.75
$ws.Range("J38").Value = 1991
$ws.Range("L38").Value = 5973
$ws.Range("N38").Value = -6717
$ws.Range("H68").Value = 30000
$ws.Range("J68").Value = 30000
$ws.Range("L68").Value = 30000
$ws.Range("N68").Value = -31498
$ws.Range("H71").Value = 30000
$ws.Range("J71").Value = 30000
$ws.Range("L71").Value = 90000
$ws.Range("N71").Value = -97488
$ws.Range("H98").Value = 606.4286
$ws.Range("I98").Value = 622.3077
$ws.Range("J98").Value = 400
$ws.Range("K98").Value = 622.3077
$ws.Range("L98").Value = 400
$ws.Range("M98").Value = 875.6923
$ws.Range("N98").Value = -3396
$ws.Range("H122").Value = 606.4286
$ws.Range("I122").Value = 622.3077
$ws.Range("J122").Value = 400
$ws.Range("K122").Value = 1866.9231
$ws.Range("L122").Value = 1200
$ws.Range("M122").Value = 583.0769
$ws.Range("N122").Value = -6100
$ws.Range("H129").Value = 1540.2916
$ws.Range("J129").Value = 1731.9048
$ws.Range("L129").Value = 5195.7144
$ws.Range("N129").Value = -15195.7144
$ws.Range("H134").Value = 47997.5
$ws.Range("J134").Value = 47997.5
$ws.Range("L134").Value = 47997.5
$ws.Range("N134").Value = -58137.5
$ws.Range("H137").Value = 65047.4
$ws.Range("I137").Value = 79554.50999999999
$ws.Range("J137").Value = 3392.1667
$ws.Range("K137").Value = 238663.53
$ws.Range("L137").Value = 10176.5001
$ws.Range("M137").Value = -236113.53
$ws.Range("N137").Value = -15276.5001
$ws.Range("H138").Value = 123210.37
$ws.Range("J138").Value = 141432.5
$ws.Range("L138").Value = 424297.5
$ws.Range("N138").Value = -434577.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8249.255999999999
$ws.Range("I32").Value = 5691.311
$ws.Range("K32").Value = 5691.311
$ws.Range("M32").Value = -5404.311
$ws.Range("H45").Value = 2530.88
$ws.Range("I45").Value = 2699.2307
$ws.Range("K45").Value = 2699.2307
$ws.Range("M45").Value = -2322.2307
$ws.Range("H61").Value = 1966.2683
$ws.Range("I61").Value = 1657.3143
$ws.Range("J61").Value = 3768.5
$ws.Range("K61").Value = 1657.3143
$ws.Range("L61").Value = 3768.5
$ws.Range("M61").Value = -1445.3143
$ws.Range("N61").Value = -4192.5
$ws.Range("H122").Value = 2421.5
$ws.Range("I122").Value = 2296.15
$ws.Range("K122").Value = 6888.450000000001
$ws.Range("M122").Value = -4438.450000000001
$ws.Range("H132").Value = 13626.182
$ws.Range("I132").Value = 2040.2667
$ws.Range("K132").Value = 6120.800099999999
$ws.Range("M132").Value = -3590.800099999999
$ws.Range("H136").Value = 1966.2683
$ws.Range("I136").Value = 1657.3143
$ws.Range("J136").Value = 3768.5
$ws.Range("K136").Value = 4971.9429
$ws.Range("L136").Value = 11305.5
$ws.Range("M136").Value = -2421.9429
$ws.Range("N136").Value = -16405.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2585.6785
$ws.Range("I20").Value = 2199.9412
$ws.Range("K20").Value = 2199.9412
$ws.Range("M20").Value = -1952.9412
$ws.Range("H107").Value = 2004.7407
$ws.Range("I107").Value = 1622.762
$ws.Range("J107").Value = 3341.6667
$ws.Range("K107").Value = 1622.762
$ws.Range("L107").Value = 3341.6667
$ws.Range("M107").Value = 297.2380000000001
$ws.Range("N107").Value = -7181.6667
$ws.Range("H134").Value = 4526.3335
$ws.Range("I134").Value = 4736.724
$ws.Range("K134").Value = 14210.172
$ws.Range("M134").Value = -11675.172

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1423.5555
$ws.Range("I16").Value = 1152
$ws.Range("J16").Value = 1966.6666
$ws.Range("K16").Value = 1152
$ws.Range("L16").Value = 1966.6666
$ws.Range("M16").Value = -865
$ws.Range("N16").Value = -2540.6666
$ws.Range("H31").Value = 3584.0544
$ws.Range("I31").Value = 1647.4062
$ws.Range("J31").Value = 6278.522
$ws.Range("K31").Value = 1647.4062
$ws.Range("L31").Value = 6278.522
$ws.Range("M31").Value = -1352.4062
$ws.Range("N31").Value = -6868.522
$ws.Range("H34").Value = 3584.0544
$ws.Range("I34").Value = 1647.4062
$ws.Range("J34").Value = 6278.522
$ws.Range("K34").Value = 1647.4062
$ws.Range("L34").Value = 6278.522
$ws.Range("M34").Value = -1445.4062
$ws.Range("N34").Value = -6682.522
$ws.Range("H58").Value = 15632.172
$ws.Range("I58").Value = 1319.25
$ws.Range("K58").Value = 1319.25
$ws.Range("M58").Value = -1116.25
$ws.Range("H113").Value = 1423.5555
$ws.Range("I113").Value = 1152
$ws.Range("J113").Value = 1966.6666
$ws.Range("K113").Value = 1152
$ws.Range("L113").Value = 1966.6666
$ws.Range("M113").Value = 1018
$ws.Range("N113").Value = -6306.6666
$ws.Range("H132").Value = 2588.7
$ws.Range("I132").Value = 1886.762
$ws.Range("J132").Value = 4226.5557
$ws.Range("K132").Value = 5660.286
$ws.Range("L132").Value = 12679.6671
$ws.Range("M132").Value = -3130.286
$ws.Range("N132").Value = -17739.6671
$ws.Range("H134").Value = 1079.1482
$ws.Range("I134").Value = 855.0833
$ws.Range("J134").Value = 2871.6667
$ws.Range("K134").Value = 2565.2499
$ws.Range("L134").Value = 8615.000100000001
$ws.Range("M134").Value = -30.2498999999998
$ws.Range("N134").Value = -13685.0001
$ws.Range("H136").Value = 15632.172
$ws.Range("I136").Value = 1319.25
$ws.Range("K136").Value = 3957.75
$ws.Range("M136").Value = -1407.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H94").Value = 3313.7896
$ws.Range("I94").Value = 1169.25
$ws.Range("J94").Value = 4873.4546
$ws.Range("K94").Value = 3507.75
$ws.Range("L94").Value = 14620.3638
$ws.Range("M94").Value = -2831.75
$ws.Range("N94").Value = -15972.3638
$ws.Range("H131").Value = 650.75
$ws.Range("I131").Value = 312.23077
$ws.Range("J131").Value = 769.6892
$ws.Range("K131").Value = 936.69231
$ws.Range("L131").Value = 2309.0676
$ws.Range("M131").Value = 4103.30769
$ws.Range("N131").Value = -12389.0676

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 18521892
$ws.Range("I80").Value = 29414742
$ws.Range("J80").Value = 4045
$ws.Range("K80").Value = 29414742
$ws.Range("L80").Value = 4045
$ws.Range("M80").Value = -29413744
$ws.Range("N80").Value = -6041
$ws.Range("H83").Value = 18521892
$ws.Range("I83").Value = 29414742
$ws.Range("J83").Value = 4045
$ws.Range("K83").Value = 147073710
$ws.Range("L83").Value = 20225
$ws.Range("M83").Value = -147068718
$ws.Range("N83").Value = -30209
$ws.Range("H113").Value = 8783.416999999999
$ws.Range("I113").Value = 12087.625
$ws.Range("J113").Value = 2175
$ws.Range("K113").Value = 12087.625
$ws.Range("L113").Value = 2175
$ws.Range("M113").Value = -9917.625
$ws.Range("N113").Value = -6515
$ws.Range("H126").Value = 5509.0454
$ws.Range("J126").Value = 4969.9
$ws.Range("L126").Value = 14909.7
$ws.Range("N126").Value = -19849.7
$ws.Range("H132").Value = 17111.422
$ws.Range("I132").Value = 5257.1763
$ws.Range("J132").Value = 26707.715
$ws.Range("K132").Value = 15771.5289
$ws.Range("L132").Value = 80123.145
$ws.Range("M132").Value = -13241.5289
$ws.Range("N132").Value = -85183.145

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 6250.5
$ws.Range("I22").Value = 10001
$ws.Range("J22").Value = 2500
$ws.Range("K22").Value = 10001
$ws.Range("L22").Value = 2500
$ws.Range("M22").Value = -9706
$ws.Range("N22").Value = -3090
$ws.Range("H27").Value = 6250.5
$ws.Range("I27").Value = 10001
$ws.Range("J27").Value = 2500
$ws.Range("K27").Value = 10001
$ws.Range("L27").Value = 2500
$ws.Range("M27").Value = -9894
$ws.Range("N27").Value = -2714
$ws.Range("H46").Value = 1950
$ws.Range("I46").Value = 400
$ws.Range("J46").Value = 3500
$ws.Range("K46").Value = 400
$ws.Range("L46").Value = 3500
$ws.Range("M46").Value = -212
$ws.Range("N46").Value = -3876
$ws.Range("H100").Value = 1697.8
$ws.Range("I100").Value = 1182
$ws.Range("K100").Value = 1182
$ws.Range("M100").Value = -641
$ws.Range("H132").Value = 209862.03
$ws.Range("I132").Value = 281809.16
$ws.Range("K132").Value = 845427.48
$ws.Range("M132").Value = -842897.48
$ws.Range("H136").Value = 2072.0588
$ws.Range("I136").Value = 1904.6875
$ws.Range("J136").Value = 4750
$ws.Range("K136").Value = 5714.0625
$ws.Range("L136").Value = 14250
$ws.Range("M136").Value = -3164.0625
$ws.Range("N136").Value = -19350

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1018
$ws.Range("I81").Value = 1096.6666
$ws.Range("J81").Value = 900
$ws.Range("K81").Value = 2193.3332
$ws.Range("L81").Value = 1800
$ws.Range("M81").Value = -1132.3332
$ws.Range("N81").Value = -3922
$ws.Range("H84").Value = 1018
$ws.Range("I84").Value = 1096.6666
$ws.Range("J84").Value = 900
$ws.Range("K84").Value = 10966.666
$ws.Range("L84").Value = 9000
$ws.Range("M84").Value = -5662.666000000001
$ws.Range("N84").Value = -19608
$ws.Range("H113").Value = 1067.15
$ws.Range("J113").Value = 592.75
$ws.Range("L113").Value = 1778.25
$ws.Range("N113").Value = -6118.25
$ws.Range("H132").Value = 658.881
$ws.Range("I132").Value = 540.5172
$ws.Range("K132").Value = 1621.5516
$ws.Range("M132").Value = 908.4484
